$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171482563018799
$ws.Range("B1").Value = 2.438419342041016
$ws.Range("D1").Value = 2.36396336555481
$ws.Range("E1").Value = 1.23833179473877
